# Lab01 rubric: rescale max points from 50 -> 40, and tidy up a few
# stray blank formatted cells left over on both the "Rubric" and "Grade"
# sheets. (commit: "notes revisions and rubric changed to 40 max")

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Rubric sheet
# ---------------------------------------------------------------
$rubric = $wb.Worksheets.Item("Rubric")

# Row 2 was just a stray blank formatted cell under the title - remove it.
$rubric.Range("A2").Style = "Normal"
$rubric.Range("A2").ClearContents()

# New point values (50-point scale -> 40-point scale).
$rubric.Range("B4").Style = "Normal"
$rubric.Range("B4").Value = 8

$rubric.Range("B5").Style = "Normal"
$rubric.Range("B5").Value = 3

$rubric.Range("B6").Style = "Normal"
$rubric.Range("B6").Value = 2

$rubric.Range("B7").Style = "Normal"
$rubric.Range("B7").Value = 2

$rubric.Range("B8").Style = "Normal"
$rubric.Range("B8").Value = 12

$rubric.Range("B9").Style = "Normal"
$rubric.Range("B9").Value = 5
$rubric.Range("C9").Style = "Normal"
$rubric.Range("C9").ClearContents()

$rubric.Range("B10").Style = "Normal"
$rubric.Range("B10").Value = 5
$rubric.Range("C10").Style = "Normal"
$rubric.Range("C10").ClearContents()

$rubric.Range("B11").Style = "Normal"
$rubric.Range("B11").Value = 3
$rubric.Range("C11").Style = "Normal"
$rubric.Range("C11").ClearContents()

# Stray blank formatted cells below the table - remove them (C12/C13 keep
# their border-ish styling, only B12/B13 are cleared).
$rubric.Range("B12").Style = "Normal"
$rubric.Range("B12").ClearContents()

$rubric.Range("B13").Style = "Normal"
$rubric.Range("B13").ClearContents()

# Total row: make the trailing C14 cell match the bold-italic "Total"
# emphasis already used by A14/B14 (was plain italic).
$rubric.Range("C14").Font.Bold = $true
$rubric.Range("C14").Font.Italic = $true

# ---------------------------------------------------------------
# Grade sheet
# ---------------------------------------------------------------
$grade = $wb.Worksheets.Item("Grade")

$grade.Range("A2").Style = "Normal"

$grade.Range("A4").Style = "Normal"

# Row 5 was just a stray blank formatted cell - remove it.
$grade.Range("A5").Style = "Normal"
$grade.Range("A5").ClearContents()

$grade.Range("B7").Style = "Normal"
$grade.Range("B7").Value = 8
$grade.Range("C7").Style = "Normal"
$grade.Range("C7").Value = 8

$grade.Range("B8").Style = "Normal"
$grade.Range("B8").Value = 3
$grade.Range("C8").Style = "Normal"
$grade.Range("C8").Value = 3

$grade.Range("B9").Style = "Normal"
$grade.Range("B9").Value = 2
$grade.Range("C9").Style = "Normal"
$grade.Range("C9").Value = 2

$grade.Range("B10").Style = "Normal"
$grade.Range("B10").Value = 2
$grade.Range("C10").Style = "Normal"
$grade.Range("C10").Value = 2

$grade.Range("B11").Style = "Normal"
$grade.Range("B11").Value = 12
$grade.Range("C11").Style = "Normal"
$grade.Range("C11").Value = 12

$grade.Range("B12").Style = "Normal"
$grade.Range("B12").Value = 5
$grade.Range("C12").Style = "Normal"
$grade.Range("C12").Value = 5

$grade.Range("B13").Style = "Normal"
$grade.Range("B13").Value = 5
$grade.Range("C13").Style = "Normal"
$grade.Range("C13").Value = 5

$grade.Range("B14").Style = "Normal"
$grade.Range("B14").Value = 3
$grade.Range("C14").Style = "Normal"
$grade.Range("C14").Value = 3

# Stray blank formatted cells - remove them.
$grade.Range("B15").Style = "Normal"
$grade.Range("B15").ClearContents()
$grade.Range("C15").Style = "Normal"
$grade.Range("C15").ClearContents()

$grade.Range("B16").Style = "Normal"
$grade.Range("B16").ClearContents()
$grade.Range("C16").Style = "Normal"
$grade.Range("C16").ClearContents()

# New column D formatting, added alongside the cleanup.
$grade.Range("D1").EntireColumn.ColumnWidth = 17.5

# Give the Grade sheet a plain portrait page setup (as added in the diff).
$grade.PageSetup.Orientation = 1

# ---------------------------------------------------------------
# Selection / active sheet (Rubric becomes the active tab)
# ---------------------------------------------------------------
$grade.Range("E16").Select()
$rubric.Select()
$rubric.Range("A1:C14").Select()
